$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Header cells Q1 ("drop") and R1 ("dropExplanation") ---
# Build their look by cloning the existing header format (P1) and then
# adjusting the fill so it matches the new fill definition used by the
# source workbook (solid FF99CCFF foreground / FF000000 background).
$p1 = $ws2.Range("P1")
$q1 = $ws2.Range("Q1")
$r1 = $ws2.Range("R1")

$p1.Copy()
$q1.PasteSpecial(-4122)   # xlPasteFormats
$q1.Interior.PatternColor = 0
$q1.Interior.Color = 16764057   # RGB(153,204,255) == FF99CCFF

$p1.Copy()
$r1.PasteSpecial(-4122)   # xlPasteFormats
$r1.Interior.PatternColor = 0
$r1.Interior.Color = 16764057   # RGB(153,204,255) == FF99CCFF
$r1.Borders.Item(7).LineStyle = 0    # clear left border inherited from P1
$r1.Borders.Item(10).Color = 0
$r1.Borders.Item(10).LineStyle = 1   # thin black right border only

$q1.Value = "drop"
$r1.Value = "dropExplanation"

# --- Body cells Q2:R41 ---
# Establish the plain-black-Calibri body style used for the drop columns by
# cloning an untouched default-styled cell, then recoloring its font.
$blank = $ws2.Range("Z99")
$blank.Copy()
$body = $ws2.Range("Q2:R41")
$body.PasteSpecial(-4122)   # xlPasteFormats
$body.Font.Color = 0

$dropCol = $ws2.Range("Q2:Q41")
$dropCol.Value = $false

# --- Selection / active sheet bookkeeping ---
$ws2.Activate()
$ws2.Range("Q1:R41").Select()
